$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts old N->O, old P->Q)
$ws.Columns("N").Insert()

# The newly inserted column inherits the column width of column M
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab and update its selection
$ws.Activate() | Out-Null
$ws.Range("S6").Select() | Out-Null
